# Regenerate the "K" column (column G) values for save_data.
# These are the strikeout counts (K) computed per-row; the sheet previously
# held a different "Strike#" style value in this column. Write the new
# K values row by row (rows 2-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 2
    4  = 3
    5  = 0
    6  = 3
    7  = 1
    8  = 1
    9  = 1
    10 = 6
    11 = 1
    12 = 3
    13 = 1
    14 = 4
    15 = 3
    16 = 5
    17 = 7
    18 = 6
    19 = 3
    20 = 2
    21 = 3
    22 = 4
    23 = 3
    24 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
